$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 37870.133
$ws.Range("I98").Value = 47166.832
$ws.Range("J98").Value = 683.3333
$ws.Range("K98").Value = 47166.832
$ws.Range("L98").Value = 683.3333
$ws.Range("M98").Value = -45668.832
$ws.Range("N98").Value = -3679.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 37870.133
$ws.Range("I122").Value = 47166.832
$ws.Range("J122").Value = 683.3333
$ws.Range("K122").Value = 141500.496
$ws.Range("L122").Value = 2049.9999
$ws.Range("M122").Value = -139050.496
$ws.Range("N122").Value = -6949.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2269325.2
$ws.Range("I132").Value = 2647231.8
$ws.Range("J132").Value = 1886.2222
$ws.Range("K132").Value = 7941695.399999999
$ws.Range("L132").Value = 5658.6666
$ws.Range("M132").Value = -7939165.399999999
$ws.Range("N132").Value = -10718.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1451.262
$ws.Range("I137").Value = 712.3043
$ws.Range("J137").Value = 2345.7896
$ws.Range("K137").Value = 2136.9129
$ws.Range("L137").Value = 7037.3688
$ws.Range("M137").Value = 413.0870999999997
$ws.Range("N137").Value = -12137.3688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 40000810
$ws.Range("I2").Value = 58824176
$ws.Range("J2").Value = 1152.875
$ws.Range("K2").Value = 58824176
$ws.Range("L2").Value = 1152.875
$ws.Range("M2").Value = -58824063
$ws.Range("N2").Value = -1378.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1783.68
$ws.Range("I61").Value = 1352.3158
$ws.Range("J61").Value = 3149.6667
$ws.Range("K61").Value = 1352.3158
$ws.Range("L61").Value = 3149.6667
$ws.Range("M61").Value = -1140.3158
$ws.Range("N61").Value = -3573.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 40000810
$ws.Range("I116").Value = 58824176
$ws.Range("J116").Value = 1152.875
$ws.Range("K116").Value = 58824176
$ws.Range("L116").Value = 1152.875
$ws.Range("M116").Value = -58821882
$ws.Range("N116").Value = -5740.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1783.68
$ws.Range("I136").Value = 1352.3158
$ws.Range("J136").Value = 3149.6667
$ws.Range("K136").Value = 4056.9474
$ws.Range("L136").Value = 9449.000100000001
$ws.Range("M136").Value = -1506.9474
$ws.Range("N136").Value = -14549.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 40000810
$ws.Range("I3").Value = 58824176
$ws.Range("J3").Value = 1152.875
$ws.Range("K3").Value = 58824176
$ws.Range("L3").Value = 1152.875
$ws.Range("M3").Value = -58824062
$ws.Range("N3").Value = -1380.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1130
$ws.Range("I99").Value = 1151.6666
$ws.Range("K99").Value = 1151.6666
$ws.Range("M99").Value = 346.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2876650.5
$ws.Range("I31").Value = 2392.5278
$ws.Range("J31").Value = 7579981.5
$ws.Range("K31").Value = 2392.5278
$ws.Range("L31").Value = 7579981.5
$ws.Range("M31").Value = -2097.5278
$ws.Range("N31").Value = -7580571.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2876650.5
$ws.Range("I34").Value = 2392.5278
$ws.Range("J34").Value = 7579981.5
$ws.Range("K34").Value = 2392.5278
$ws.Range("L34").Value = 7579981.5
$ws.Range("M34").Value = -2190.5278
$ws.Range("N34").Value = -7580385.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 43333.332
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 43333.332
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 712.5
$ws.Range("I107").Value = 684.8276
$ws.Range("J107").Value = 759.7059
$ws.Range("K107").Value = 684.8276
$ws.Range("L107").Value = 759.7059
$ws.Range("M107").Value = 1235.1724
$ws.Range("N107").Value = -4599.7059

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 867.03845
$ws.Range("I122").Value = 804.4211
$ws.Range("J122").Value = 1037
$ws.Range("K122").Value = 2413.2633
$ws.Range("L122").Value = 3111
$ws.Range("M122").Value = 36.73669999999993
$ws.Range("N122").Value = -8011

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 620.2963
$ws.Range("J5").Value = 1164.091
$ws.Range("L5").Value = 3492.273
$ws.Range("N5").Value = -3716.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 680
$ws.Range("J41").Value = 750
$ws.Range("L41").Value = 2250
$ws.Range("N41").Value = -2926

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 428.75
$ws.Range("I92").Value = 220
$ws.Range("J92").Value = 637.5
$ws.Range("K92").Value = 660
$ws.Range("L92").Value = 1912.5
$ws.Range("M92").Value = 588
$ws.Range("N92").Value = -4408.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 288.66666
$ws.Range("J97").Value = 319.6
$ws.Range("L97").Value = 958.8000000000001
$ws.Range("N97").Value = -1950.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 620.2963
$ws.Range("J135").Value = 1164.091
$ws.Range("L135").Value = 10476.819
$ws.Range("N135").Value = -15546.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 42131850
$ws.Range("I137").Value = 25650624
$ws.Range("J137").Value = 52844644
$ws.Range("K137").Value = 76951872
$ws.Range("L137").Value = 158533932
$ws.Range("M137").Value = -76946772
$ws.Range("N137").Value = -158544132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3756
$ws.Range("I80").Value = 4438.5713
$ws.Range("J80").Value = 2163.3333
$ws.Range("K80").Value = 4438.5713
$ws.Range("L80").Value = 2163.3333
$ws.Range("M80").Value = -3440.5713
$ws.Range("N80").Value = -4159.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3756
$ws.Range("I83").Value = 4438.5713
$ws.Range("J83").Value = 2163.3333
$ws.Range("K83").Value = 22192.8565
$ws.Range("L83").Value = 10816.6665
$ws.Range("M83").Value = -17200.8565
$ws.Range("N83").Value = -20800.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 229.5
$ws.Range("I55").Value = 215.26315
$ws.Range("J55").Value = 247.53334
$ws.Range("K55").Value = 215.26315
$ws.Range("L55").Value = 247.53334
$ws.Range("M55").Value = -42.26315
$ws.Range("N55").Value = -593.53334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3537.85
$ws.Range("I81").Value = 2845.125
$ws.Range("J81").Value = 3999.6667
$ws.Range("K81").Value = 5690.25
$ws.Range("L81").Value = 7999.3334
$ws.Range("M81").Value = -4629.25
$ws.Range("N81").Value = -10121.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3537.85
$ws.Range("I84").Value = 2845.125
$ws.Range("J84").Value = 3999.6667
$ws.Range("K84").Value = 28451.25
$ws.Range("L84").Value = 39996.667
$ws.Range("M84").Value = -23147.25
$ws.Range("N84").Value = -50604.667
